# RoA Master Sheet - refresh T12M RoA (col D) and T6M RoA (col E) figures
# with the latest pull, plus a couple of small ITD RoA (col C) rounding
# corrections. Also reset the sheet view back to the top-left and leave
# the selection on the first empty row below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - AIRCRAFT F1
$ws.Cells.Item(2,4).Value = 0.23129999999999998
$ws.Cells.Item(2,5).Value = 0.21840000000000001

# Row 6 - 1.0 LEGACY ABS F1
$ws.Cells.Item(6,4).Value = 0.34619999999999995
$ws.Cells.Item(6,5).Value = 0.32890000000000003

# Row 7 - 1L EETC F1
$ws.Cells.Item(7,3).Value = 0.14529999999999998
$ws.Cells.Item(7,4).Value = 0.07980000000000001
$ws.Cells.Item(7,5).Value = 0.0926

# Row 8 - 2L EETC F1
$ws.Cells.Item(8,4).Value = 0.1804
$ws.Cells.Item(8,5).Value = 0.14940000000000001

# Row 9 - 3.0 MEZZ ABS F1
$ws.Cells.Item(9,3).Value = 0.045599999999999995
$ws.Cells.Item(9,4).Value = 0.045599999999999995
$ws.Cells.Item(9,5).Value = 0.045599999999999995

# Row 10 - 3.0 SENIOR ABS F1
$ws.Cells.Item(10,4).Value = 0.1454
$ws.Cells.Item(10,5).Value = 0.19500000000000001

# Row 11 - CMBS F1 (substrategy subtotal)
$ws.Cells.Item(11,4).Value = 0.07580000000000001
$ws.Cells.Item(11,5).Value = 0.0736

# Row 12 - AIR UNSECURED F1
$ws.Cells.Item(12,4).Value = 0.048499999999999995
$ws.Cells.Item(12,5).Value = 0.049

# Row 13 - AIRCRAFT F1_INCOME
$ws.Cells.Item(13,4).Value = 0.0867
$ws.Cells.Item(13,5).Value = 0.1067

# Row 14 - CMBS 2.0/3.0 IG F1
$ws.Cells.Item(14,4).Value = -0.0675
$ws.Cells.Item(14,5).Value = -0.0508

# Row 15 - CMBS 2.0/3.0 NON-IG F1
$ws.Cells.Item(15,3).Value = 0.039599999999999996
$ws.Cells.Item(15,4).Value = 0.019799999999999998
$ws.Cells.Item(15,5).Value = 0.029300000000000003

# Row 16 - CMBS AGENCY F1
$ws.Cells.Item(16,4).Value = 0.0102
$ws.Cells.Item(16,5).Value = 0.009899999999999999

# Row 17 - CMBS IO F1
$ws.Cells.Item(17,3).Value = 0.16440000000000002
$ws.Cells.Item(17,4).Value = 0.0165
$ws.Cells.Item(17,5).Value = 0.013500000000000002

# Row 18 - CMBS PRIVATE LOANS
$ws.Cells.Item(18,4).Value = 0.0144
$ws.Cells.Item(18,5).Value = 0.0069

# Row 19 - SHORT TERM (substrategy subtotal)
$ws.Cells.Item(19,4).Value = 0.0046
$ws.Cells.Item(19,5).Value = 0.0039000000000000003

# Row 20 - CMBS SASB F1
$ws.Cells.Item(20,4).Value = 0.0046
$ws.Cells.Item(20,5).Value = 0.0039000000000000003

# Row 21 - CLO F1 (strategy subtotal)
$ws.Cells.Item(21,4).Value = 0.029500000000000002
$ws.Cells.Item(21,5).Value = 0.029500000000000002

# Row 22 - CLO AAA ETF F1
$ws.Cells.Item(22,4).Value = 0.029500000000000002
$ws.Cells.Item(22,5).Value = 0.029500000000000002

# Row 24 - ABS F1
$ws.Cells.Item(24,4).Value = 0.0137
$ws.Cells.Item(24,5).Value = 0.0137

# Row 26 - SENIOR MPL
$ws.Cells.Item(26,3).Value = 0.026000000000000002

# Scroll the sheet view back to the top and leave the cursor on D27,
# just below the data table.
$ws.Range("D27").Select()
